$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.784.13'
$ws.Range('E2').Value = '  +0.95%  '
$ws.Range('D3').Value = '1.864.81'
$ws.Range('E3').Value = '  +0.91%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.023'
$ws.Range('E4').Value = '  -1.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '321.94'
$ws.Range('E5').Value = '  +0.30%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.022'
$ws.Range('E6').Value = '  -0.66%  '
$ws.Range('E7').Value = '  +0.28%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3820'
$ws.Range('E8').Value = '  +1.53%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07481'
$ws.Range('E9').Value = '  +1.35%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8849'
$ws.Range('E10').Value = '  +1.44%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '21.65'
$ws.Range('E11').Value = '  +1.33%  '
$ws.Range('D12').Value = '1.861.58'
$ws.Range('E12').Value = '  +0.47%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.771'
$ws.Range('E13').Value = '  +1.59%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.511'
$ws.Range('E14').Value = '  +0.21%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.07128'
$ws.Range('E15').Value = '  -0.75%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '86.61'
$ws.Range('E16').Value = '  +5.05%  '
$ws.Range('E17').Value = '  -0.73%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000009098'
$ws.Range('E18').Value = '  +0.68%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.022'
$ws.Range('E19').Value = '  -0.68%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '15.52'
$ws.Range('E20').Value = '  +0.96%  '
$ws.Range('D21').Value = '27.782.95'
$ws.Range('E21').Value = '  +0.89%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.300'
$ws.Range('E22').Value = '  +1.27%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.20'
$ws.Range('E23').Value = '  -1.00%  '
$ws.Range('D24').Value = '2.090.55'
$ws.Range('E24').Value = '  +0.78%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.057'
$ws.Range('E25').Value = '  +7.04%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '158.07'
$ws.Range('E26').Value = '  +0.39%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.78'
$ws.Range('E27').Value = '  +0.76%  '
$ws.Range('B28').Value = 'LidoDAOToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.012'
$ws.Range('E28').Value = '  +3.20%  '
$ws.Range('B29').Value = 'InternetComputer(DFINITY)'
$ws.Range('C29').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.386'
$ws.Range('E29').Value = '  +2.54%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '121.92'
$ws.Range('E30').Value = '  +4.55%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.09072'
$ws.Range('E31').Value = '  +0.56%  '
$ws.Range('E32').Value = '  +1.93%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.7692'
$ws.Range('E33').Value = '  +1.32%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.051'
$ws.Range('E34').Value = '  +6.11%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.573'
$ws.Range('E35').Value = '  +1.83%  '
$ws.Range('E36').Value = '  -0.60%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.146'
$ws.Range('E37').Value = '  -0.23%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01985'
$ws.Range('E38').Value = '  +0.76%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.05304'
$ws.Range('E39').Value = '  +0.56%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.876'
$ws.Range('E40').Value = '  +2.65%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.5212'
$ws.Range('E41').Value = '  +1.51%  '
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.952'
$ws.Range('E42').Value = '  +3.93%  '
$ws.Range('B43').Value = 'Algorand'
$ws.Range('C43').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1685'
$ws.Range('E43').Value = '  +1.05%  '
$ws.Range('B44').Value = 'Aptos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.726'
$ws.Range('E44').Value = '  +3.31%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '10.81'
$ws.Range('E45').Value = '  +2.31%  '
$ws.Range('B46').Value = 'Quant'
$ws.Range('C46').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '110.30'
$ws.Range('E46').Value = '  +1.35%  '
$ws.Range('B47').Value = 'NEARProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.720'
$ws.Range('E47').Value = '  +1.12%  '
$ws.Range('B48').Value = 'PaxDollar'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.025'
$ws.Range('E48').Value = '  -0.60%  '
$ws.Range('B49').Value = 'Decentraland'
$ws.Range('C49').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.4728'
$ws.Range('E49').Value = '  +2.20%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06507'
$ws.Range('E50').Value = '  +1.64%  '
$ws.Range('B51').Value = 'RenderToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.887'
$ws.Range('E51').Value = '  +1.90%  '
